$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-23 Sunday" "2025-03-24 Monday"
Replace-Text "542×7=" "902×8="
Replace-Text "695×4=" "683×8="
Replace-Text "681×3=" "560×2="
Replace-Text "341×7=" "965×6="
Replace-Text "229×3=" "442×9="
Replace-Text "742×9=" "246×5="
Replace-Text "826×4=" "511×2="
Replace-Text "739×3=" "531×2="
Replace-Text "464×4=" "437×8="
Replace-Text "878×2=" "982×3="
Replace-Text "158×7=" "713×5="
Replace-Text "506×2=" "122×2="
Replace-Text "976×5=" "823×4="
Replace-Text "401×6=" "526×9="
Replace-Text "399×4=" "107×7="
Replace-Text "375×3=" "639×8="
Replace-Text "822×7=" "708×2="
Replace-Text "253×7=" "743×4="
Replace-Text "914×6=" "873×4="
Replace-Text "362×5=" "878×4="
Replace-Text "342×4=" "423×2="
Replace-Text "428×4=" "530×9="
Replace-Text "217×8=" "775×7="
Replace-Text "621×5=" "435×5="
Replace-Text "538×4=" "973×7="
